$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("pop").Name = "population"
$wb.Worksheets.Item("pop_births_deaths").Name = "population_births_deaths"
$wb.Worksheets.Item("pop_missing_axis_name").Name = "population_missing_axis_name"
$wb.Worksheets.Item("pop_missing_values").Name = "population_missing_values"
$wb.Worksheets.Item("pop_narrow_format").Name = "population_narrow_format"
